$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 46082 = 2026-03-01) for every
# data row (rows 2 through 411). The update bumps that date forward by one day
# (serial 46083 = 2026-03-02) for all of them.
$ws.Range("C2:C411").Value = 46083
